$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.397.20'
$ws.Range("E2").Value = '  +4.34%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.487.72'
$ws.Range("E3").Value = '  +3.64%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.06'
$ws.Range("E5").Value = '  +2.89%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.18'
$ws.Range("E6").Value = '  +7.74%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("E8").Value = '  +1.53%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.71'
$ws.Range("E9").Value = '  +0.22%  '

# Row 10
$ws.Range("E10").Value = '  +4.76%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.084.21'
$ws.Range("E12").Value = '  +3.71%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.81'
$ws.Range("E13").Value = '  +7.04%  '

# Row 14
$ws.Range("E14").Value = '  -0.10%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.493.02'
$ws.Range("E15").Value = '  +3.77%  '

# Row 16
$ws.Range("E16").Value = '  +4.13%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.388.24'
$ws.Range("E17").Value = '  +4.07%  '

# Row 18
$ws.Range("E18").Value = '  +3.19%  '

# Row 19
$ws.Range("E19").Value = '  +6.42%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.41'
$ws.Range("E20").Value = '  +5.78%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '393.40'
$ws.Range("E21").Value = '  +2.89%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.566'
$ws.Range("E22").Value = '  +3.42%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.26'
$ws.Range("E23").Value = '  +0.06%  '

# Row 24
$ws.Range("E24").Value = '  -0.04%  '

# Row 25
$ws.Range("E25").Value = '  +8.77%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.628.71'
$ws.Range("E26").Value = '  +3.69%  '

# Row 27
$ws.Range("E27").Value = '  -2.75%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.84'
$ws.Range("E28").Value = '  +10.22%  '

# Row 29
$ws.Range("E29").Value = '  -0.10%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.29'
$ws.Range("E30").Value = '  +5.75%  '

# Row 31
$ws.Range("E31").Value = '  +2.45%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.45'
$ws.Range("E32").Value = '  +9.02%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.00%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.83'
$ws.Range("E34").Value = '  +3.90%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '32.54'
$ws.Range("E35").Value = '  +27.67%  '

# Row 36
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.34'
$ws.Range("E36").Value = '  +8.84%  '

# Row 37
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.16'
$ws.Range("E37").Value = '  +5.25%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '171.26'
$ws.Range("E38").Value = '  +2.48%  '

# Row 39
$ws.Range("E39").Value = '  +9.99%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.525.84'
$ws.Range("E40").Value = '  +3.72%  '

# Row 41
$ws.Range("E41").Value = '  +1.86%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.805'
$ws.Range("E42").Value = '  +4.70%  '

# Row 43
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.51'
$ws.Range("E43").Value = '  +4.26%  '

# Row 44
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.74'
$ws.Range("E44").Value = '  +7.64%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.54'
$ws.Range("E45").Value = '  +0.75%  '

# Row 46
$ws.Range("E46").Value = '  +10.50%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.619.39'
$ws.Range("E47").Value = '  +7.48%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.91'
$ws.Range("E48").Value = '  +8.36%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.31'
$ws.Range("E49").Value = '  +18.99%  '

# Row 50
$ws.Range("E50").Value = '  +2.18%  '

# Row 51
$ws.Range("E51").Value = '  +5.38%  '
